# Auto-generated Excel COM-interop script
# Applies cell updates for commit: 'Add data for 2022-08-17'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 4468
$ws.Range('I3').Value = 4673
$ws.Range('D4').Value = 1930
$ws.Range('I4').Value = 1079
$ws.Range('I5').Value = 432
$ws.Range('I6').Value = 5088
$ws.Range('D7').Value = 28120
$ws.Range('I7').Value = 15740

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I2').Value = 127
$ws.Range('I8').Value = 960
$ws.Range('I12').Value = 36
$ws.Range('I14').Value = 85
$ws.Range('I15').Value = 181
$ws.Range('I18').Value = 109
$ws.Range('I19').Value = 444
$ws.Range('I21').Value = 83
$ws.Range('I25').Value = 77
$ws.Range('I27').Value = 147
$ws.Range('I29').Value = 998
$ws.Range('I33').Value = 729
$ws.Range('I34').Value = 75
$ws.Range('I37').Value = 507
$ws.Range('I42').Value = 536
$ws.Range('I44').Value = 115
$ws.Range('I48').Value = 217
$ws.Range('I52').Value = 335
$ws.Range('I55').Value = 172
$ws.Range('I58').Value = 11
$ws.Range('I60').Value = 81
$ws.Range('D63').Value = 320
$ws.Range('I63').Value = 68
$ws.Range('I64').Value = 139
$ws.Range('I65').Value = 355
$ws.Range('I67').Value = 619
$ws.Range('I68').Value = 56
$ws.Range('I72').Value = 58
$ws.Range('I73').Value = 133
$ws.Range('I76').Value = 237
$ws.Range('I77').Value = 95
$ws.Range('I78').Value = 223
$ws.Range('I79').Value = 440
$ws.Range('I80').Value = 54
$ws.Range('I83').Value = 320
$ws.Range('I85').Value = 704
$ws.Range('I86').Value = 93
$ws.Range('I89').Value = 180
$ws.Range('I95').Value = 261
$ws.Range('I97').Value = 121
$ws.Range('I99').Value = 296
$ws.Range('D101').Value = 28120
$ws.Range('I101').Value = 15740

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 181
$ws.Range('I3').Value = 282
$ws.Range('I6').Value = 178
$ws.Range('I7').Value = 704

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I3').Value = 119
$ws.Range('I5').Value = 12
$ws.Range('I7').Value = 335

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 301
$ws.Range('I3').Value = 273
$ws.Range('I4').Value = 57
$ws.Range('I6').Value = 303
$ws.Range('I7').Value = 960

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I6').Value = 65
$ws.Range('I7').Value = 180

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I2').Value = 32
$ws.Range('I7').Value = 85

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I3').Value = 162
$ws.Range('I7').Value = 507

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 80
$ws.Range('I3').Value = 108
$ws.Range('I7').Value = 296

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 146
$ws.Range('I4').Value = 38
$ws.Range('I7').Value = 619

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I5').Value = 17
$ws.Range('I7').Value = 355

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I2').Value = 110
$ws.Range('I3').Value = 124
$ws.Range('I7').Value = 320

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I3').Value = 100
$ws.Range('I7').Value = 261

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I2').Value = 165
$ws.Range('I6').Value = 231
$ws.Range('I7').Value = 729

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I3').Value = 72
$ws.Range('I4').Value = 24

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 289
$ws.Range('I3').Value = 347
$ws.Range('I7').Value = 998

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I2').Value = 165
$ws.Range('I6').Value = 121
$ws.Range('I7').Value = 444

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I3').Value = 32
$ws.Range('I7').Value = 115

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I6').Value = 124
$ws.Range('I7').Value = 217

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I2').Value = 50
$ws.Range('I7').Value = 237

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 139
$ws.Range('I3').Value = 184
$ws.Range('I6').Value = 151
$ws.Range('I7').Value = 536

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I6').Value = 84
$ws.Range('I7').Value = 223

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I3').Value = 54
$ws.Range('I7').Value = 172

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('I6').Value = 60
$ws.Range('I7').Value = 83

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I3').Value = 140
$ws.Range('I6').Value = 128
$ws.Range('I7').Value = 440

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('I2').Value = 39
$ws.Range('I6').Value = 48
$ws.Range('I7').Value = 139

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('I2').Value = 32
$ws.Range('I7').Value = 109

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('I2').Value = 31
$ws.Range('I7').Value = 75

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('I3').Value = 23
$ws.Range('I7').Value = 77

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I3').Value = 42
$ws.Range('I5').Value = 7
$ws.Range('I7').Value = 181

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I6').Value = 32
$ws.Range('I7').Value = 133

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('I6').Value = 27
$ws.Range('I7').Value = 127

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I3').Value = 23
$ws.Range('I7').Value = 121

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('I6').Value = 58
$ws.Range('I7').Value = 147

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I4').Value = 45
$ws.Range('I7').Value = 93

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('I2').Value = 21
$ws.Range('I7').Value = 56

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I2').Value = 27
$ws.Range('I6').Value = 23
$ws.Range('I7').Value = 81

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('I2').Value = 10
$ws.Range('I7').Value = 58

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('I2').Value = 30
$ws.Range('I3').Value = 33
$ws.Range('I7').Value = 95

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('I3').Value = 11
$ws.Range('I7').Value = 54

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('I2').Value = 8
$ws.Range('I7').Value = 36

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range('I6').Value = 5
$ws.Range('I7').Value = 11
